# Update countries & provincias Spain
#
# This script applies the daily COVID-19 data refresh to the "Pais" sheet:
#  - Updates case totals (columns B-H) for several countries whose numbers changed.
#  - Because the country ranking shifted with the new totals, a few rows also
#    need their country label (column A) updated so that the label stays in
#    sync with the correct numbers for that rank position:
#      * Niger overtakes Libano (rows 99-101 shuffle: Niger, Libano, Costa Rica)
#      * Burundi overtakes Nicaragua (rows 199-204 shuffle down by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 1158881
$ws.Range("C4").Value = 27851
$ws.Range("D4").Value = 160668
$ws.Range("E4").Value = 930920
$ws.Range("F4").Value = 16475
$ws.Range("G4").Value = 1540
$ws.Range("H4").Value = 67293

# Row 15: Canada -> Canada
$ws.Range("B15").Value = 56714
$ws.Range("C15").Value = 1653
$ws.Range("D15").Value = 23621
$ws.Range("E15").Value = 29529
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 173
$ws.Range("H15").Value = 3564

# Row 45: Noruega -> Noruega
$ws.Range("B45").Value = 7809
$ws.Range("C45").Value = 26
$ws.Range("D45").Value = 32
$ws.Range("E45").Value = 7566
$ws.Range("F45").Value = 37
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 211

# Row 46: Chequia -> Chequia
$ws.Range("B46").Value = 7755
$ws.Range("C46").Value = 18
$ws.Range("D46").Value = 3461
$ws.Range("E46").Value = 4049
$ws.Range("F46").Value = 67
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 245

# Row 99: Libano -> Niger
$ws.Range("A99").Value = "Niger"
$ws.Range("B99").Value = 736
$ws.Range("C99").Value = 8
$ws.Range("D99").Value = 507
$ws.Range("E99").Value = 194
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 35

# Row 100: Costa Rica -> Libano
$ws.Range("A100").Value = "Libano"
$ws.Range("B100").Value = 733
$ws.Range("C100").Value = 4
$ws.Range("D100").Value = 197
$ws.Range("E100").Value = 511
$ws.Range("F100").Value = 43
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 25

# Row 101: Niger -> Costa Rica
$ws.Range("A101").Value = "Costa Rica"
$ws.Range("B101").Value = 733
$ws.Range("C101").Value = 8
$ws.Range("D101").Value = 372
$ws.Range("E101").Value = 355
$ws.Range("F101").Value = 6
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 6

# Row 122: Reunion -> Reunion
$ws.Range("B122").Value = 423
$ws.Range("C122").Value = 1
$ws.Range("D122").Value = 300
$ws.Range("E122").Value = 123
$ws.Range("F122").Value = 2
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

# Row 168: Libia -> Libia
$ws.Range("B168").Value = 63
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 22
$ws.Range("E168").Value = 38
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 3

# Row 199: Nicaragua -> Burundi
$ws.Range("A199").Value = "Burundi"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 4
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

# Row 200: Islas Malvinas -> Nicaragua
$ws.Range("A200").Value = "Nicaragua"
$ws.Range("B200").Value = 14
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 7
$ws.Range("E200").Value = 4
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 3

# Row 201: Islas Turcas y Caicos -> Islas Malvinas
$ws.Range("A201").Value = "Islas Malvinas"
$ws.Range("B201").Value = 13
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 13
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# Row 202: Santa Sede -> Islas Turcas y Caicos
$ws.Range("A202").Value = "Islas Turcas y Caicos"
$ws.Range("B202").Value = 12
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 5
$ws.Range("E202").Value = 6
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 1

# Row 203: Montserrat -> Santa Sede
$ws.Range("A203").Value = "Santa Sede"
$ws.Range("B203").Value = 11
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204: Burundi -> Montserrat
$ws.Range("A204").Value = "Montserrat"
$ws.Range("B204").Value = 11
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 2
$ws.Range("E204").Value = 8
$ws.Range("F204").Value = 1
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 1

